$wb = $excel.ActiveWorkbook

# Rename "Spouse 2" sheet to "Casey"
$spouseSheet = $wb.Worksheets.Item("Spouse 2")
$spouseSheet.Name = "Casey"

# Make "Casey" the active/selected sheet (moves tabSelected + workbookView.activeTab)
$spouseSheet.Activate()
